$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.081747055053711
$ws.Range("B1").Value = 2.411656141281128
$ws.Range("C1").Value = 5.037034034729004
$ws.Range("D1").Value = 2.289585113525391
$ws.Range("E1").Value = 1.295222640037537
